# -----------------------------------------------------------------------
# Applies two changes described by the diff:
#   1. Remove the stray "_GoBack" bookmark that sits after the
#      "Faltan puntajes en perfil de usuario." paragraph.
#   2. Expand "Mostrar más datos (puntaje de la publicación)" into a run
#      split that highlights "puntaje de la publicación" in yellow, adds
#      ", puntaje del cliente" and moves the "_GoBack" bookmark to sit
#      right before the closing parenthesis.
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Remove the old "_GoBack" bookmark -------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Rework "Mostrar más datos (puntaje de la publicación)" ---------
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$found = $findRange.Find.Execute(
    "Mostrar más datos (puntaje de la publicación)",
    $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)

if ($found) {
    $startPos = $findRange.Start

    # Remove the original single run's text; we'll rebuild it piece by piece.
    $findRange.Text = ""

    # Piece 1 (plain): "Mostrar más datos ("
    $r1 = $d.Range($startPos, $startPos)
    $r1.InsertAfter("Mostrar más datos (")
    $pos = $r1.End

    # Piece 2 (highlighted yellow): "puntaje de la publicación"
    $r2 = $d.Range($pos, $pos)
    $r2.InsertAfter("puntaje de la publicación")
    $r2 = $d.Range($pos, $r2.End)
    $r2.Font.HighlightColorIndex = "#FFFF00"
    $pos = $r2.End

    # Piece 3 (plain): ", puntaje del cliente"
    $r3 = $d.Range($pos, $pos)
    $r3.InsertAfter(", puntaje del cliente")
    $pos = $r3.End

    # Piece 4 (plain): ")" -- inserted before the bookmark so that the
    # bookmark's collapsed range no longer sits immediately before the
    # paragraph mark (doing so first avoids an anchor-resolution quirk).
    $r4 = $d.Range($pos, $pos)
    $r4.InsertAfter(")")

    # Re-create the "_GoBack" bookmark right between ", puntaje del cliente"
    # and ")".
    $bmRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
